# Relocating the return_slider definition moved the variable's
# initialization earlier in the code, which changed the starting point
# (and therefore the converged result) of the portfolio optimizer.
# Update the recalculated "Opt Portfolio" (C) and "Opt Portfolio with
# View" (D) columns with the new optimizer output values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.000000000000000002927345865710862
$ws.Range("D2").Value = 0

$ws.Range("C3").Value = 0.02124717618409318
$ws.Range("D3").Value = 0.02124715427314639

$ws.Range("C4").Value = 0.1273382014015918
$ws.Range("D4").Value = 0.1273388940376063

$ws.Range("C5").Value = 0.1697704453391681
$ws.Range("D5").Value = 0.1697699223641621

$ws.Range("C6").Value = 0.2390021039874927
$ws.Range("D6").Value = 0.2390017823104368

$ws.Range("C7").Value = 0.09744142026537399
$ws.Range("D7").Value = 0.09744136897287403

$ws.Range("C8").Value = 0.34520065282228
$ws.Range("D8").Value = 0.3452008780417743
